$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "Sittz" -> "Sitz" in the cell describing player creation
$ws.Range("A5").Value = "Funktion die einen Spieler grundsätzlich erstellt und in die DB einträgt (mit Namen, Balance, an welchem Sitz er am Tisch sitzt)"

# Move the active selection to B6 (was A20)
$ws.Range("B6").Select()
